$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so formatted numeric
# strings (e.g. "1.00", "61.352.42") are preserved exactly as authored.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.352.42'
$ws.Range("E2").Value = '  -1.36%  '
$ws.Range("D3").Value = '2.985.07'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '595.30'
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").Value = '143.60'
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '2.984.61'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("D10").Value = '0.147'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  +4.88%  '
$ws.Range("E12").Value = '  +3.53%  '
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '34.25'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '0.125'
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").Value = '3.476.39'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = '61.338.60'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '6.88'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").Value = '2.984.38'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").Value = '449.27'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("E21").Value = '  +2.10%  '
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '7.31'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").Value = '81.87'
$ws.Range("E24").Value = '  +1.75%  '
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = '10.55'
$ws.Range("E25").Value = '  +4.66%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("D27").Value = '11.98'
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +3.25%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +1.98%  '
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").Value = '27.23'
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").Value = '0.0₃0821'
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '5.78'
$ws.Range("E37").Value = '  +1.50%  '
$ws.Range("D38").Value = '50.25'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").Value = '9.00'
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.120'
$ws.Range("E42").Value = '  +6.84%  '
$ws.Range("D43").Value = '386.79'
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0351'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.269'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = '38.45'
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("D47").Value = '2.694.51'
$ws.Range("E47").Value = '  -1.92%  '
$ws.Range("D48").Value = '130.71'
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '2.14'
$ws.Range("E51").Value = '  +0.88%  '
